$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '49.578.99'
$ws.Range('E2').Value = '  -0.79%  '
$ws.Range('D3').Value = '2.646.73'
$ws.Range('E3').Value = '  +0.02%  '
$ws.Range('E4').Value = '  -0.03%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '112.63'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  -0.99%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '326.66'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  -0.02%  '
$ws.Range('E7').Value = '  -1.13%  '
$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '0.999'
$c.Style = 'Normal'
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('E9').Value = '  -1.52%  '
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '39.70'
$c.Style = 'Normal'
$ws.Range('E10').Value = '  -3.33%  '
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '20.03'
$c.Style = 'Normal'
$ws.Range('E11').Value = '  -0.86%  '
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '0.0815'
$c.Style = 'Normal'
$ws.Range('E12').Value = '  -0.93%  '
$ws.Range('E13').Value = '  +1.88%  '
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '7.60'
$c.Style = 'Normal'
$ws.Range('E14').Value = '  +2.76%  '
$ws.Range('D15').Value = '3.059.79'
$ws.Range('E15').Value = '  -0.08%  '
$ws.Range('D16').Value = '2.640.95'
$ws.Range('E16').Value = '  -1.37%  '
$ws.Range('E17').Value = '  -1.60%  '
$ws.Range('D18').Value = '49.564.54'
$ws.Range('E18').Value = '  -0.67%  '
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '13.38'
$c.Style = 'Normal'
$ws.Range('E19').Value = '  +1.18%  '
$ws.Range('E20').Value = '  -1.74%  '
$ws.Range('E21').Value = '  -0.52%  '
$ws.Range('D22').Value = '0.0₃0949'
$ws.Range('E22').Value = '  -1.06%  '
$ws.Range('E23').Value = '  -3.00%  '
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '68.85'
$c.Style = 'Normal'
$ws.Range('E24').Value = '  -4.46%  '
$ws.Range('E25').Value = '  -0.68%  '
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '26.11'
$c.Style = 'Normal'
$ws.Range('E26').Value = '  -2.73%  '
$ws.Range('E27').Value = '  +0.01%  '
$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '10.14'
$c.Style = 'Normal'
$ws.Range('E28').Value = '  +1.43%  '
$ws.Range('E29').Value = '  -0.94%  '
$ws.Range('E30').Value = '  -2.71%  '
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '34.67'
$c.Style = 'Normal'
$ws.Range('E31').Value = '  -4.18%  '
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '49.63'
$c.Style = 'Normal'
$ws.Range('E32').Value = '  -1.39%  '
$ws.Range('E33').Value = '  +0.37%  '
$ws.Range('E34').Value = '  +1.59%  '
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '19.22'
$c.Style = 'Normal'
$ws.Range('E35').Value = '  -1.25%  '
$ws.Range('E36').Value = '  -0.15%  '
$ws.Range('E37').Value = '  -1.44%  '
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '2.04'
$c.Style = 'Normal'
$ws.Range('E38').Value = '  -1.57%  '
$ws.Range('E39').Value = '  +0.44%  '
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '129.89'
$c.Style = 'Normal'
$ws.Range('E40').Value = '  +5.11%  '
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '23.70'
$c.Style = 'Normal'
$ws.Range('E41').Value = '  +7.42%  '
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '0.0350'
$c.Style = 'Normal'
$ws.Range('E42').Value = '  +11.33%  '
$ws.Range('E43').Value = '  +2.89%  '
$ws.Range('E44').Value = '  -0.76%  '
$ws.Range('D45').Value = '2.063.85'
$ws.Range('E45').Value = '  -1.01%  '
$ws.Range('E46').Value = '  -0.72%  '
$ws.Range('E47').Value = '  +6.58%  '
$ws.Range('E48').Value = '  -4.87%  '
$ws.Range('E49').Value = '  -2.39%  '
$ws.Range('E50').Value = '  -2.95%  '
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '58.72'
$c.Style = 'Normal'
$ws.Range('E51').Value = '  -2.15%  '
